# Scheduled market-data refresh: updates the Universalis-derived price/profit
# columns (H:N) for the leves whose item market data changed since the last run.
# Only the specific cells that changed are touched; every other cell (including
# leve name/item/level/exp/gil/amount/item-id in columns A:G) is left untouched.
$wb = $excel.ActiveWorkbook

# ===== ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 62: "The Mustache Suits Him" (Enchanted Mythrite Ink)
$ws.Cells.Item(62, 8).Value = 4271.4  # H62 currentAveragePrice: 4384.3335 -> 4271.4
$ws.Cells.Item(62, 9).Value = 2407.1  # I62 currentAveragePriceNQ: 2576.5 -> 2407.1
$ws.Cells.Item(62, 11).Value = 2407.1  # K62 LevePriceNQ: 2576.5 -> 2407.1
$ws.Cells.Item(62, 13).Value = -1783.1  # M62 LeveProfitNQ: -1952.5 -> -1783.1

# Row 65: "Forgery of Convenience (L)" (Enchanted Mythrite Ink)
$ws.Cells.Item(65, 8).Value = 4271.4  # H65 currentAveragePrice: 4384.3335 -> 4271.4
$ws.Cells.Item(65, 9).Value = 2407.1  # I65 currentAveragePriceNQ: 2576.5 -> 2407.1
$ws.Cells.Item(65, 11).Value = 12035.5  # K65 LevePriceNQ: 12882.5 -> 12035.5
$ws.Cells.Item(65, 13).Value = -8915.5  # M65 LeveProfitNQ: -9762.5 -> -8915.5

# Row 80: "Cleansing the Wicked Humours" (Hallowed Water)
$ws.Cells.Item(80, 8).Value = 25397.666  # H80 currentAveragePrice: 19129 -> 25397.666
$ws.Cells.Item(80, 9).Value = 532.3333  # I80 currentAveragePriceNQ: 474.75 -> 532.3333
$ws.Cells.Item(80, 10).Value = 50263  # J80 currentAveragePriceHQ: 37783.25 -> 50263
$ws.Cells.Item(80, 11).Value = 1596.9999  # K80 LevePriceNQ: 1424.25 -> 1596.9999
$ws.Cells.Item(80, 12).Value = 150789  # L80 LevePriceHQ: 113349.75 -> 150789
$ws.Cells.Item(80, 13).Value = -598.9999  # M80 LeveProfitNQ: -426.25 -> -598.9999
$ws.Cells.Item(80, 14).Value = -152785  # N80 LeveProfitHQ: -115345.75 -> -152785

# Row 83: "Washing Away the Sins (L)" (Hallowed Water)
$ws.Cells.Item(83, 8).Value = 25397.666  # H83 currentAveragePrice: 19129 -> 25397.666
$ws.Cells.Item(83, 9).Value = 532.3333  # I83 currentAveragePriceNQ: 474.75 -> 532.3333
$ws.Cells.Item(83, 10).Value = 50263  # J83 currentAveragePriceHQ: 37783.25 -> 50263
$ws.Cells.Item(83, 11).Value = 4790.9997  # K83 LevePriceNQ: 4272.75 -> 4790.9997
$ws.Cells.Item(83, 12).Value = 452367  # L83 LevePriceHQ: 340049.25 -> 452367
$ws.Cells.Item(83, 13).Value = 201.0002999999997  # M83 LeveProfitNQ: 719.25 -> 201.0002999999997
$ws.Cells.Item(83, 14).Value = -462351  # N83 LeveProfitHQ: -350033.25 -> -462351

# Row 101: "Edge of the Arcane" (Cunning Craftsman's Tea)
$ws.Cells.Item(101, 8).Value = 410.5  # H101 currentAveragePrice: 372.6 -> 410.5
$ws.Cells.Item(101, 9).Value = 410.5  # I101 currentAveragePriceNQ: 372.6 -> 410.5
$ws.Cells.Item(101, 11).Value = 1231.5  # K101 LevePriceNQ: 1117.8 -> 1231.5
$ws.Cells.Item(101, 13).Value = 390.5  # M101 LeveProfitNQ: 504.1999999999998 -> 390.5

# Row 132: "Fast-forwarding Flora" (Growth Formula Lambda)
$ws.Cells.Item(132, 8).Value = 1543.5  # H132 currentAveragePrice: 1602.84 -> 1543.5
$ws.Cells.Item(132, 9).Value = 1483.591  # I132 currentAveragePriceNQ: 1606.2 -> 1483.591
$ws.Cells.Item(132, 10).Value = 1873  # J132 currentAveragePriceHQ: 1589.4 -> 1873
$ws.Cells.Item(132, 11).Value = 4450.772999999999  # K132 LevePriceNQ: 4818.6 -> 4450.772999999999
$ws.Cells.Item(132, 12).Value = 5619  # L132 LevePriceHQ: 4768.200000000001 -> 5619
$ws.Cells.Item(132, 13).Value = -1920.772999999999  # M132 LeveProfitNQ: -2288.6 -> -1920.772999999999
$ws.Cells.Item(132, 14).Value = -10679  # N132 LeveProfitHQ: -9828.200000000001 -> -10679

# Row 137: "Cutting Edge of Culinary Quality" (Magnesia Whetstone)
$ws.Cells.Item(137, 8).Value = 2492.7273  # H137 currentAveragePrice: 3008.2354 -> 2492.7273
$ws.Cells.Item(137, 9).Value = 1146  # I137 currentAveragePriceNQ: 1436 -> 1146
$ws.Cells.Item(137, 11).Value = 3438  # K137 LevePriceNQ: 4308 -> 3438
$ws.Cells.Item(137, 13).Value = -888  # M137 LeveProfitNQ: -1758 -> -888

# Row 138: "All-night Crafting" (Cunning Craftsman's Tisane)
$ws.Cells.Item(138, 8).Value = 5848.9766  # H138 currentAveragePrice: 6755.0356 -> 5848.9766
$ws.Cells.Item(138, 9).Value = 4232.6665  # I138 currentAveragePriceNQ: 4057.6 -> 4232.6665
$ws.Cells.Item(138, 10).Value = 5970.2  # J138 currentAveragePriceHQ: 7341.4346 -> 5970.2
$ws.Cells.Item(138, 11).Value = 12697.9995  # K138 LevePriceNQ: 12172.8 -> 12697.9995
$ws.Cells.Item(138, 12).Value = 17910.6  # L138 LevePriceHQ: 22024.3038 -> 17910.6
$ws.Cells.Item(138, 13).Value = -7557.999500000002  # M138 LeveProfitNQ: -7032.799999999999 -> -7557.999500000002
$ws.Cells.Item(138, 14).Value = -28190.6  # N138 LeveProfitHQ: -32304.3038 -> -28190.6

# ===== ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 32: "Ingot We Trust" (Steel Ingot)
$ws.Cells.Item(32, 8).Value = 21138.127  # H32 currentAveragePrice: 21979 -> 21138.127
$ws.Cells.Item(32, 9).Value = 14112.934  # I32 currentAveragePriceNQ: 15942.846 -> 14112.934
$ws.Cells.Item(32, 11).Value = 14112.934  # K32 LevePriceNQ: 15942.846 -> 14112.934
$ws.Cells.Item(32, 13).Value = -13825.934  # M32 LeveProfitNQ: -15655.846 -> -13825.934

# Row 61: "Dealing with the Tough Stuff" (Cobalt Ingot)
$ws.Cells.Item(61, 8).Value = 1934.7826  # H61 currentAveragePrice: 2004.3478 -> 1934.7826
$ws.Cells.Item(61, 9).Value = 1638.8889  # I61 currentAveragePriceNQ: 1652.9412 -> 1638.8889
$ws.Cells.Item(61, 11).Value = 1638.8889  # K61 LevePriceNQ: 1652.9412 -> 1638.8889
$ws.Cells.Item(61, 13).Value = -1426.8889  # M61 LeveProfitNQ: -1440.9412 -> -1426.8889

# Row 122: "Haste for High Durium" (High Durium Nugget)
$ws.Cells.Item(122, 8).Value = 558570.0600000001  # H122 currentAveragePrice: 669904.0600000001 -> 558570.0600000001
$ws.Cells.Item(122, 9).Value = 1002026.1  # I122 currentAveragePriceNQ: 1430651.6 -> 1002026.1
$ws.Cells.Item(122, 11).Value = 3006078.3  # K122 LevePriceNQ: 4291954.800000001 -> 3006078.3
$ws.Cells.Item(122, 13).Value = -3003628.3  # M122 LeveProfitNQ: -4289504.800000001 -> -3003628.3

# Row 132: "Don't Bore Me, Ore Me" (Mountain Chromite Ingot)
$ws.Cells.Item(132, 8).Value = 1803.7812  # H132 currentAveragePrice: 1770.7576 -> 1803.7812
$ws.Cells.Item(132, 10).Value = 2987.3333  # J132 currentAveragePriceHQ: 2419 -> 2987.3333
$ws.Cells.Item(132, 12).Value = 8961.999899999999  # L132 LevePriceHQ: 7257 -> 8961.999899999999
$ws.Cells.Item(132, 14).Value = -14021.9999  # N132 LeveProfitHQ: -12317 -> -14021.9999

# Row 136: "Metal with Mettle" (Cobalt Tungsten Ingot)
$ws.Cells.Item(136, 8).Value = 1934.7826  # H136 currentAveragePrice: 2004.3478 -> 1934.7826
$ws.Cells.Item(136, 9).Value = 1638.8889  # I136 currentAveragePriceNQ: 1652.9412 -> 1638.8889
$ws.Cells.Item(136, 11).Value = 4916.6667  # K136 LevePriceNQ: 4958.8236 -> 4916.6667
$ws.Cells.Item(136, 13).Value = -2366.6667  # M136 LeveProfitNQ: -2408.8236 -> -2366.6667

# ===== BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 86: "Through Thick and Thin" (Adamantite Nugget)
$ws.Cells.Item(86, 8).Value = 3361.4  # H86 currentAveragePrice: 3300.3333 -> 3361.4
$ws.Cells.Item(86, 10).Value = 6807  # J86 currentAveragePriceHQ: 4901 -> 6807
$ws.Cells.Item(86, 12).Value = 6807  # L86 LevePriceHQ: 4901 -> 6807
$ws.Cells.Item(86, 14).Value = -9053  # N86 LeveProfitHQ: -7147 -> -9053

# Row 89: "Piercing Eyes Deserve Piercing Shafts (L)" (Adamantite Nugget)
$ws.Cells.Item(89, 8).Value = 3361.4  # H89 currentAveragePrice: 3300.3333 -> 3361.4
$ws.Cells.Item(89, 10).Value = 6807  # J89 currentAveragePriceHQ: 4901 -> 6807
$ws.Cells.Item(89, 12).Value = 34035  # L89 LevePriceHQ: 24505 -> 34035
$ws.Cells.Item(89, 14).Value = -45267  # N89 LeveProfitHQ: -35737 -> -45267

# Row 134: "Ruthenium Supremium" (Ruthenium Ingot)
$ws.Cells.Item(134, 8).Value = 1521.9286  # H134 currentAveragePrice: 1439.52 -> 1521.9286
$ws.Cells.Item(134, 9).Value = 695.8333  # I134 currentAveragePriceNQ: 699.4783 -> 695.8333
$ws.Cells.Item(134, 10).Value = 6478.5  # J134 currentAveragePriceHQ: 9950 -> 6478.5
$ws.Cells.Item(134, 11).Value = 2087.4999  # K134 LevePriceNQ: 2098.4349 -> 2087.4999
$ws.Cells.Item(134, 12).Value = 19435.5  # L134 LevePriceHQ: 29850 -> 19435.5
$ws.Cells.Item(134, 13).Value = 447.5001000000002  # M134 LeveProfitNQ: 436.5650999999998 -> 447.5001000000002
$ws.Cells.Item(134, 14).Value = -24505.5  # N134 LeveProfitHQ: -34920 -> -24505.5

# ===== CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 31: "Wall Not Found" (Walnut Lumber)
$ws.Cells.Item(31, 8).Value = 5276.4375  # H31 currentAveragePrice: 4875.8 -> 5276.4375
$ws.Cells.Item(31, 9).Value = 2141  # I31 currentAveragePriceNQ: 1896.8572 -> 2141
$ws.Cells.Item(31, 10).Value = 6321.5835  # J31 currentAveragePriceHQ: 7482.375 -> 6321.5835
$ws.Cells.Item(31, 11).Value = 2141  # K31 LevePriceNQ: 1896.8572 -> 2141
$ws.Cells.Item(31, 12).Value = 6321.5835  # L31 LevePriceHQ: 7482.375 -> 6321.5835
$ws.Cells.Item(31, 13).Value = -1846  # M31 LeveProfitNQ: -1601.8572 -> -1846
$ws.Cells.Item(31, 14).Value = -6911.5835  # N31 LeveProfitHQ: -8072.375 -> -6911.5835

# Row 34: "Armoires of the Rich and Famous" (Walnut Lumber)
$ws.Cells.Item(34, 8).Value = 5276.4375  # H34 currentAveragePrice: 4875.8 -> 5276.4375
$ws.Cells.Item(34, 9).Value = 2141  # I34 currentAveragePriceNQ: 1896.8572 -> 2141
$ws.Cells.Item(34, 10).Value = 6321.5835  # J34 currentAveragePriceHQ: 7482.375 -> 6321.5835
$ws.Cells.Item(34, 11).Value = 2141  # K34 LevePriceNQ: 1896.8572 -> 2141
$ws.Cells.Item(34, 12).Value = 6321.5835  # L34 LevePriceHQ: 7482.375 -> 6321.5835
$ws.Cells.Item(34, 13).Value = -1939  # M34 LeveProfitNQ: -1694.8572 -> -1939
$ws.Cells.Item(34, 14).Value = -6725.5835  # N34 LeveProfitHQ: -7886.375 -> -6725.5835

# Row 132: "Hull Lotta Damage" (Ginseng Lumber)
$ws.Cells.Item(132, 8).Value = 3172.5908  # H132 currentAveragePrice: 2994.2917 -> 3172.5908
$ws.Cells.Item(132, 9).Value = 2612.1428  # I132 currentAveragePriceNQ: 2442.6 -> 2612.1428
$ws.Cells.Item(132, 10).Value = 4153.375  # J132 currentAveragePriceHQ: 3913.7778 -> 4153.375
$ws.Cells.Item(132, 11).Value = 7836.428400000001  # K132 LevePriceNQ: 7327.799999999999 -> 7836.428400000001
$ws.Cells.Item(132, 12).Value = 12460.125  # L132 LevePriceHQ: 11741.3334 -> 12460.125
$ws.Cells.Item(132, 13).Value = -5306.428400000001  # M132 LeveProfitNQ: -4797.799999999999 -> -5306.428400000001
$ws.Cells.Item(132, 14).Value = -17520.125  # N132 LeveProfitHQ: -16801.3334 -> -17520.125

# Row 134: "Wood You Be Quiet" (Ceiba Lumber)
$ws.Cells.Item(134, 8).Value = 3228.9285  # H134 currentAveragePrice: 2535.65 -> 3228.9285
$ws.Cells.Item(134, 9).Value = 1671.3334  # I134 currentAveragePriceNQ: 1407.1428 -> 1671.3334
$ws.Cells.Item(134, 10).Value = 6032.6  # J134 currentAveragePriceHQ: 5168.8335 -> 6032.6
$ws.Cells.Item(134, 11).Value = 5014.0002  # K134 LevePriceNQ: 4221.428400000001 -> 5014.0002
$ws.Cells.Item(134, 12).Value = 18097.8  # L134 LevePriceHQ: 15506.5005 -> 18097.8
$ws.Cells.Item(134, 13).Value = -2479.0002  # M134 LeveProfitNQ: -1686.428400000001 -> -2479.0002
$ws.Cells.Item(134, 14).Value = -23167.8  # N134 LeveProfitHQ: -20576.5005 -> -23167.8

# ===== CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 4: "In Hot Water" (Boiled Egg)
$ws.Cells.Item(4, 8).Value = 118402110  # H4 currentAveragePrice: 76925230 -> 118402110
$ws.Cells.Item(4, 9).Value = 154831310  # I4 currentAveragePriceNQ: 111111900 -> 154831310
$ws.Cells.Item(4, 10).Value = 7193.5  # J4 currentAveragePriceHQ: 5221.75 -> 7193.5
$ws.Cells.Item(4, 11).Value = 464493930  # K4 LevePriceNQ: 333335700 -> 464493930
$ws.Cells.Item(4, 12).Value = 21580.5  # L4 LevePriceHQ: 15665.25 -> 21580.5
$ws.Cells.Item(4, 13).Value = -464493818  # M4 LeveProfitNQ: -333335588 -> -464493818
$ws.Cells.Item(4, 14).Value = -21804.5  # N4 LeveProfitHQ: -15889.25 -> -21804.5

# Row 23: "Sweet Smell of Success" (Lavender Oil)
$ws.Cells.Item(23, 8).Value = 254.66667  # H23 currentAveragePrice: 344.5 -> 254.66667
$ws.Cells.Item(23, 10).Value = 82  # J23 currentAveragePriceHQ: 89 -> 82
$ws.Cells.Item(23, 12).Value = 246  # L23 LevePriceHQ: 267 -> 246
$ws.Cells.Item(23, 14).Value = -716  # N23 LeveProfitHQ: -737 -> -716

# Row 68: "Such a Butter Face" (Fermented Butter)
$ws.Cells.Item(68, 8).Value = 1810.8182  # H68 currentAveragePrice: 1819.909 -> 1810.8182
$ws.Cells.Item(68, 9).Value = 1638  # I68 currentAveragePriceNQ: 1658 -> 1638
$ws.Cells.Item(68, 11).Value = 4914  # K68 LevePriceNQ: 4974 -> 4914
$ws.Cells.Item(68, 13).Value = -4103  # M68 LeveProfitNQ: -4163 -> -4103

# Row 71: "No Margarine of Error (L)" (Fermented Butter)
$ws.Cells.Item(71, 8).Value = 1810.8182  # H71 currentAveragePrice: 1819.909 -> 1810.8182
$ws.Cells.Item(71, 9).Value = 1638  # I71 currentAveragePriceNQ: 1658 -> 1638
$ws.Cells.Item(71, 11).Value = 14742  # K71 LevePriceNQ: 14922 -> 14742
$ws.Cells.Item(71, 13).Value = -10686  # M71 LeveProfitNQ: -10866 -> -10686

# ===== GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 105: "Untucked" (Palladium Tuck)
$ws.Cells.Item(105, 8).Value = 44057  # H105 currentAveragePrice: 20814.334 -> 44057
$ws.Cells.Item(105, 10).Value = 44057  # J105 currentAveragePriceHQ: 20814.334 -> 44057
$ws.Cells.Item(105, 12).Value = 44057  # L105 LevePriceHQ: 20814.334 -> 44057
$ws.Cells.Item(105, 14).Value = -51045  # N105 LeveProfitHQ: -27802.334 -> -51045

# Row 132: "On Board for Lar" (Lar Ingot)
$ws.Cells.Item(132, 8).Value = 1511.5  # H132 currentAveragePrice: 1454.6666 -> 1511.5
$ws.Cells.Item(132, 9).Value = 1537  # I132 currentAveragePriceNQ: 1469.875 -> 1537
$ws.Cells.Item(132, 11).Value = 4611  # K132 LevePriceNQ: 4409.625 -> 4611
$ws.Cells.Item(132, 13).Value = -2081  # M132 LeveProfitNQ: -1879.625 -> -2081

# ===== LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 16: "Saddle Sore" (Hard Leather)
$ws.Cells.Item(16, 8).Value = 5838.2  # H16 currentAveragePrice: 6523.3 -> 5838.2
$ws.Cells.Item(16, 9).Value = 4438.067  # I16 currentAveragePriceNQ: 5351.533 -> 4438.067
$ws.Cells.Item(16, 11).Value = 4438.067  # K16 LevePriceNQ: 5351.533 -> 4438.067
$ws.Cells.Item(16, 13).Value = -4268.067  # M16 LeveProfitNQ: -5181.533 -> -4268.067

# Row 40: "Best Served Toad" (Toad Leather)
$ws.Cells.Item(40, 8).Value = 2811.5  # H40 currentAveragePrice: 2936.5 -> 2811.5
$ws.Cells.Item(40, 9).Value = 2915.3333  # I40 currentAveragePriceNQ: 2936.5 -> 2915.3333
$ws.Cells.Item(40, 10).Value = 2500  # J40 currentAveragePriceHQ: 0 -> 2500
$ws.Cells.Item(40, 11).Value = 2915.3333  # K40 LevePriceNQ: 2936.5 -> 2915.3333
$ws.Cells.Item(40, 12).Value = 2500  # L40 LevePriceHQ: 0 -> 2500
$ws.Cells.Item(40, 13).Value = -2779.3333  # M40 LeveProfitNQ: -2800.5 -> -2779.3333
$ws.Cells.Item(40, 14).Value = -2772  # N40 LeveProfitHQ: (blank) -> -2772

# Row 132: "Tenets of Tanning" (Silver Lobo Leather)
$ws.Cells.Item(132, 8).Value = 5310.778  # H132 currentAveragePrice: 5444.846 -> 5310.778
$ws.Cells.Item(132, 9).Value = 5310.778  # I132 currentAveragePriceNQ: 5254.091 -> 5310.778
$ws.Cells.Item(132, 10).Value = 0  # J132 currentAveragePriceHQ: 6494 -> 0
$ws.Cells.Item(132, 11).Value = 15932.334  # K132 LevePriceNQ: 15762.273 -> 15932.334
$ws.Cells.Item(132, 12).Value = 0  # L132 LevePriceHQ: 19482 -> 0
$ws.Cells.Item(132, 13).Value = -13402.334  # M132 LeveProfitNQ: -13232.273 -> -13402.334
$ws.Cells.Item(132, 14).ClearContents()  # N132 LeveProfitHQ: -24542 -> (cleared)

# ===== WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 92: "Modest Beginnings" (Bloodhempen Culottes of Casting)
$ws.Cells.Item(92, 8).Value = 200000  # H92 currentAveragePrice: 88333.336 -> 200000
$ws.Cells.Item(92, 10).Value = 200000  # J92 currentAveragePriceHQ: 88333.336 -> 200000
$ws.Cells.Item(92, 12).Value = 200000  # L92 LevePriceHQ: 88333.336 -> 200000
$ws.Cells.Item(92, 14).Value = -204992  # N92 LeveProfitHQ: -93325.336 -> -204992

# Row 113: "A Tender Table" (Pixie Floss)
$ws.Cells.Item(113, 8).Value = 1841.1428  # H113 currentAveragePrice: 1926.32 -> 1841.1428
$ws.Cells.Item(113, 9).Value = 1636.8572  # I113 currentAveragePriceNQ: 1726.75 -> 1636.8572
$ws.Cells.Item(113, 10).Value = 2045.4286  # J113 currentAveragePriceHQ: 2110.5386 -> 2045.4286
$ws.Cells.Item(113, 11).Value = 4910.571599999999  # K113 LevePriceNQ: 5180.25 -> 4910.571599999999
$ws.Cells.Item(113, 12).Value = 6136.2858  # L113 LevePriceHQ: 6331.6158 -> 6136.2858
$ws.Cells.Item(113, 13).Value = -2740.571599999999  # M113 LeveProfitNQ: -3010.25 -> -2740.571599999999
$ws.Cells.Item(113, 14).Value = -10476.2858  # N113 LeveProfitHQ: -10671.6158 -> -10476.2858

# Row 132: "Comfy Cabins" (Snow Cotton Cloth)
$ws.Cells.Item(132, 8).Value = 1517.3158  # H132 currentAveragePrice: 1568 -> 1517.3158
$ws.Cells.Item(132, 10).Value = 605  # J132 currentAveragePriceHQ: 0 -> 605
$ws.Cells.Item(132, 12).Value = 1815  # L132 LevePriceHQ: 0 -> 1815
$ws.Cells.Item(132, 14).Value = -6875  # N132 LeveProfitHQ: (blank) -> -6875
